$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value would otherwise be auto-converted from text
# into a numeric value by Excel (e.g. "42.64") must be pre-formatted as
# text so the literal string content is preserved, matching the source data.
foreach ($cellRef in @('D5', 'D8', 'D13', 'D14', 'D18', 'D20', 'D21', 'D24', 'D25', 'D27', 'D33', 'D36', 'D40', 'D42', 'D43', 'D46', 'D47', 'D48')) {
    $ws.Range($cellRef).NumberFormat = '@'
}

$ws.Range('D2').Value = '35.288.73'
$ws.Range('E2').Value = '  +1.09%  '
$ws.Range('D3').Value = '1.865.23'
$ws.Range('E4').Value = '  +0.62%  '
$ws.Range('D5').Value = '239.61'
$ws.Range('E5').Value = '  +3.35%  '
$ws.Range('E6').Value = '  +0.68%  '
$ws.Range('E7').Value = '  +0.57%  '
$ws.Range('D8').Value = '42.64'
$ws.Range('E8').Value = '  +7.02%  '
$ws.Range('E9').Value = '  +0.91%  '
$ws.Range('E10').Value = '  +1.11%  '
$ws.Range('E11').Value = '  +0.37%  '
$ws.Range('D12').Value = '2.134.41'
$ws.Range('E12').Value = '  +1.30%  '
$ws.Range('D13').Value = '11.54'
$ws.Range('E13').Value = '  +0.93%  '
$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D14').Value = '0.680'
$ws.Range('E14').Value = '  +1.06%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '1.845.29'
$ws.Range('E15').Value = '  +0.28%  '
$ws.Range('E16').Value = '  +1.80%  '
$ws.Range('D17').Value = '35.279.03'
$ws.Range('E17').Value = '  +1.06%  '
$ws.Range('D18').Value = '70.14'
$ws.Range('E18').Value = '  +0.34%  '
$ws.Range('D20').Value = '241.30'
$ws.Range('E20').Value = '  +0.28%  '
$ws.Range('D21').Value = '12.28'
$ws.Range('E21').Value = '  +0.86%  '
$ws.Range('E22').Value = '  +1.38%  '
$ws.Range('E23').Value = '  +0.58%  '
$ws.Range('D24').Value = '2.25'
$ws.Range('E24').Value = '  -1.45%  '
$ws.Range('D25').Value = '169.64'
$ws.Range('E25').Value = '  -1.06%  '
$ws.Range('E26').Value = '  +25.39%  '
$ws.Range('D27').Value = '8.12'
$ws.Range('E27').Value = '  +4.31%  '
$ws.Range('E28').Value = '  +1.72%  '
$ws.Range('E30').Value = '  +1.84%  '
$ws.Range('E31').Value = '  +0.58%  '
$ws.Range('E32').Value = '  +2.16%  '
$ws.Range('D33').Value = '1.83'
$ws.Range('E33').Value = '  +27.49%  '
$ws.Range('E34').Value = '  +2.24%  '
$ws.Range('E35').Value = '  +8.87%  '
$ws.Range('D36').Value = '0.816'
$ws.Range('E36').Value = '  +17.08%  '
$ws.Range('E37').Value = '  +5.97%  '
$ws.Range('E38').Value = '  +3.98%  '
$ws.Range('E39').Value = '  +4.43%  '
$ws.Range('D40').Value = '90.47'
$ws.Range('E40').Value = '  -0.70%  '
$ws.Range('D41').Value = '1.345.74'
$ws.Range('E41').Value = '  +0.23%  '
$ws.Range('D42').Value = '15.27'
$ws.Range('E42').Value = '  +3.14%  '
$ws.Range('D43').Value = '0.0602'
$ws.Range('E43').Value = '  +15.16%  '
$ws.Range('E44').Value = '  +2.50%  '
$ws.Range('E45').Value = '  +0.72%  '
$ws.Range('D46').Value = '12.45'
$ws.Range('E46').Value = '  +45.91%  '
$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D47').Value = '6.63'
$ws.Range('E47').Value = '  +5.15%  '
$ws.Range('B48').Value = 'MXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D48').Value = '2.74'
$ws.Range('E48').Value = '  -0.92%  '
$ws.Range('D49').Value = '2.051.19'
$ws.Range('E49').Value = '  +1.47%  '
$ws.Range('E50').Value = '  +3.34%  '
$ws.Range('E51').Value = '  +1.05%  '
